# The original "gross photo synthesis" label (rows 2-11) was a typo/formatting
# mistake - missing the correct spacing. Correct the label text, which also
# causes the shared-strings table to be de-duplicated/reordered around the
# remaining labels ("net photosynthesis", "dark respiration").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = "gross photosynthesis"
}

# Update the active selection to match the author's final cursor position.
$ws.Range("B2:B11").Select()
